$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name
$ws.Name = "Through 2021-12-11"

# Update December row label
$ws.Range("A13").Value = "December (through 12-11)"

# Update November row, 2021 column (H)
$ws.Range("H11").Value = 197

# Update December row values
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = 39
$ws.Range("F13").Value = 16
$ws.Range("G13").Value = 55
$ws.Range("H13").Value = 84

# Update Total row values
$ws.Range("B14").Value = 301
$ws.Range("C14").Value = 593
$ws.Range("D14").Value = 860
$ws.Range("F14").Value = 550
$ws.Range("G14").Value = 1319
$ws.Range("H14").Value = 1729
